# Refresh the live coin-market snapshot (price / 1h volume, plus two
# swapped rank positions) pulled in by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '58.167.62'
$ws.Range('E2').Value = '  -3.13%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.961.63'
$ws.Range('E3').Value = '  +0.66%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.40%  '

# Row 5: BNB
$ws.Range('D5').Value = '''556.38'
$ws.Range('E5').Value = '  -2.09%  '

# Row 6: Solana
$ws.Range('D6').Value = '''130.82'
$ws.Range('E6').Value = '  +7.25%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.35%  '

# Row 8: XRP
$ws.Range('D8').Value = '''0.518'
$ws.Range('E8').Value = '  +4.63%  '

# Row 9: LidoStakedEther
$ws.Range('D9').Value = '2.947.26'
$ws.Range('E9').Value = '  -0.05%  '

# Row 10: Dogecoin
$ws.Range('D10').Value = '''0.129'
$ws.Range('E10').Value = '  -1.30%  '

# Row 11: Toncoin
$ws.Range('D11').Value = '''4.81'
$ws.Range('E11').Value = '  -4.93%  '

# Row 12: Cardano
$ws.Range('D12').Value = '''0.447'
$ws.Range('E12').Value = '  +3.39%  '

# Row 13: ShibaInu
$ws.Range('D13').Value = '''0.0000223'
$ws.Range('E13').Value = '  +1.63%  '

# Row 14: Avalanche
$ws.Range('D14').Value = '''32.75'
$ws.Range('E14').Value = '  +1.84%  '

# Row 15: TRON
$ws.Range('E15').Value = '  +2.16%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Range('D16').Value = '3.452.47'
$ws.Range('E16').Value = '  +0.68%  '

# Row 17: Polkadot
$ws.Range('D17').Value = '''6.76'
$ws.Range('E17').Value = '  +10.84%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '2.965.75'
$ws.Range('E18').Value = '  +0.69%  '

# Row 19: WrappedBTC
$ws.Range('D19').Value = '58.328.58'
$ws.Range('E19').Value = '  -2.82%  '

# Row 20: BitcoinCash
$ws.Range('D20').Value = '''418.15'
$ws.Range('E20').Value = '  -1.81%  '

# Row 21: Chainlink
$ws.Range('D21').Value = '''13.11'
$ws.Range('E21').Value = '  +1.53%  '

# Row 22: Polygon
$ws.Range('D22').Value = '''0.683'
$ws.Range('E22').Value = '  +4.57%  '

# Row 23: Uniswap
$ws.Range('D23').Value = '''6.95'
$ws.Range('E23').Value = '  +0.43%  '

# Row 24: InternetComputer(DFINITY)
$ws.Range('D24').Value = '''12.99'
$ws.Range('E24').Value = '  +1.96%  '

# Row 25: Litecoin
$ws.Range('D25').Value = '''79.38'
$ws.Range('E25').Value = '  +1.81%  '

# Row 26: Dai
$ws.Range('E26').Value = '  +0.15%  '

# Row 27: FirstDigitalUSD
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  +0.30%  '

# Row 28: PancakeSwap
$ws.Range('D28').Value = '''2.50'
$ws.Range('E28').Value = '  +0.55%  '

# Row 29: RenderToken
$ws.Range('D29').Value = '''7.52'
$ws.Range('E29').Value = '  +6.41%  '

# Row 30: ImmutableX
$ws.Range('D30').Value = '''1.99'
$ws.Range('E30').Value = '  +7.36%  '

# Row 31: NEARProtocol
$ws.Range('D31').Value = '''6.17'
$ws.Range('E31').Value = '  +2.29%  '

# Row 32: Hedera
$ws.Range('D32').Value = '''0.104'
$ws.Range('E32').Value = '  +13.69%  '

# Row 33: EthereumClassic
$ws.Range('D33').Value = '''25.07'
$ws.Range('E33').Value = '  +0.27%  '

# Row 34: Filecoin
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''5.64'
$ws.Range('E34').Value = '  +2.49%  '

# Row 35: Stacks
$ws.Range('B35').Value = 'Stacks'
$ws.Range('C35').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D35').Value = '''2.12'
$ws.Range('E35').Value = '  -1.57%  '

# Row 36: Mantle
$ws.Range('D36').Value = '''0.943'
$ws.Range('E36').Value = '  +0.58%  '

# Row 37: OKB
$ws.Range('D37').Value = '''48.32'
$ws.Range('E37').Value = '  -2.05%  '

# Row 38: PEPE
$ws.Range('D38').Value = '0.0₃0680'
$ws.Range('E38').Value = '  +6.16%  '

# Row 39: Cosmos
$ws.Range('D39').Value = '''8.39'
$ws.Range('E39').Value = '  +7.72%  '

# Row 40: dogwifhat
$ws.Range('D40').Value = '''2.58'
$ws.Range('E40').Value = '  +9.36%  '

# Row 41: VeChain
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.0351'
$ws.Range('E41').Value = '  -0.25%  '

# Row 42: Kaspa
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '''0.109'
$ws.Range('E42').Value = '  +1.09%  '

# Row 43: Bittensor
$ws.Range('D43').Value = '''379.57'
$ws.Range('E43').Value = '  +1.92%  '

# Row 44: Maker
$ws.Range('D44').Value = '2.649.58'
$ws.Range('E44').Value = '  +1.80%  '

# Row 46: TheGraph
$ws.Range('D46').Value = '''0.239'
$ws.Range('E46').Value = '  +3.04%  '

# Row 47: Monero
$ws.Range('D47').Value = '''121.21'
$ws.Range('E47').Value = '  +1.95%  '

# Row 48: Stellar
$ws.Range('D48').Value = '''0.109'
$ws.Range('E48').Value = '  +3.92%  '

# Row 49: Fetch.AI
$ws.Range('D49').Value = '''1.98'
$ws.Range('E49').Value = '  +2.11%  '

# Row 50: InjectiveProtocol
$ws.Range('D50').Value = '''23.39'
$ws.Range('E50').Value = '  +1.65%  '

# Row 51: ThetaToken
$ws.Range('D51').Value = '''2.00'
$ws.Range('E51').Value = '  +2.42%  '
